$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 9).Value = 0.233114009085334
$ws.Cells.Item(2, 10).Value = 0.2331140090853341
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.1698666666666667
$ws.Cells.Item(2, 14).Value = 0.5096000000000001
$ws.Cells.Item(2, 15).Value = 0.0442365680687153
$ws.Cells.Item(2, 16).Value = 0.0442365680687153
$ws.Cells.Item(2, 17).Value = 0.08282285324444447
$ws.Cells.Item(2, 18).Value = 0.7454056792000001
$ws.Cells.Item(2, 19).Value = 0.0103121637306745
$ws.Cells.Item(2, 20).Value = 0.0103121637306745
$ws.Cells.Item(3, 9).Value = 0.233114009085334
$ws.Cells.Item(3, 10).Value = 0.2331140090853341
$ws.Cells.Item(3, 15).Value = 0.09956205276773258
$ws.Cells.Item(3, 16).Value = 0.09956205276773258
$ws.Cells.Item(3, 19).Value = 0.02320930927345172
$ws.Cells.Item(3, 20).Value = 0.02320930927345172
$ws.Cells.Item(4, 9).Value = 0.233114009085334
$ws.Cells.Item(4, 10).Value = 0.2331140090853341
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 0.6666666666666666
$ws.Cells.Item(4, 13).Value = 0.655462
$ws.Cells.Item(4, 14).Value = 1.966386
$ws.Cells.Item(4, 15).Value = 0.1706949924222307
$ws.Cells.Item(4, 16).Value = 0.1706949924222307
$ws.Cells.Item(4, 17).Value = 0.3195873216246667
$ws.Cells.Item(4, 18).Value = 2.876285894622
$ws.Cells.Item(4, 19).Value = 0.03979139401433692
$ws.Cells.Item(4, 20).Value = 0.03979139401433692
$ws.Cells.Item(5, 9).Value = 0.233114009085334
$ws.Cells.Item(5, 10).Value = 0.2331140090853341
$ws.Cells.Item(5, 13).Value = 1.828665666666667
$ws.Cells.Item(5, 14).Value = 5.485997
$ws.Cells.Item(5, 15).Value = 0.476219936646915
$ws.Cells.Item(5, 16).Value = 0.476219936646915
$ws.Cells.Item(5, 17).Value = 0.8916128815354445
$ws.Cells.Item(5, 18).Value = 8.024515933819002
$ws.Cells.Item(5, 19).Value = 0.1110135386381261
$ws.Cells.Item(5, 20).Value = 0.1110135386381262
$ws.Cells.Item(6, 9).Value = 0.233114009085334
$ws.Cells.Item(6, 10).Value = 0.2331140090853341
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.07852066666666667
$ws.Cells.Item(6, 14).Value = 0.235562
$ws.Cells.Item(6, 15).Value = 0.02044830150589229
$ws.Cells.Item(6, 16).Value = 0.02044830150589229
$ws.Cells.Item(6, 17).Value = 0.03828476639711111
$ws.Cells.Item(6, 18).Value = 0.344562897574
$ws.Cells.Item(6, 19).Value = 0.004766785543024225
$ws.Cells.Item(6, 20).Value = 0.004766785543024226
$ws.Cells.Item(7, 9).Value = 0.233114009085334
$ws.Cells.Item(7, 10).Value = 0.2331140090853341
$ws.Cells.Item(7, 13).Value = 0.7251310000000001
$ws.Cells.Item(7, 14).Value = 2.175393
$ws.Cells.Item(7, 15).Value = 0.1888381485885141
$ws.Cells.Item(7, 16).Value = 0.1888381485885141
$ws.Cells.Item(7, 17).Value = 0.3535562307456667
$ws.Cells.Item(7, 18).Value = 3.182006076711001
$ws.Cells.Item(7, 19).Value = 0.04402081788572053
$ws.Cells.Item(7, 20).Value = 0.04402081788572054
$ws.Cells.Item(8, 7).Value = 1.604000333333333
$ws.Cells.Item(8, 8).Value = 4.812001
$ws.Cells.Item(8, 9).Value = 0.7668859909146659
$ws.Cells.Item(8, 10).Value = 0.7668859909146659
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 0.1698666666666667
$ws.Cells.Item(8, 14).Value = 0.5096000000000001
$ws.Cells.Item(8, 15).Value = 0.0442365680687153
$ws.Cells.Item(8, 16).Value = 0.0442365680687153
$ws.Cells.Item(8, 17).Value = 0.2724661899555556
$ws.Cells.Item(8, 18).Value = 2.452195709600001
$ws.Cells.Item(8, 19).Value = 0.0339244043380408
$ws.Cells.Item(8, 20).Value = 0.0339244043380408
$ws.Cells.Item(9, 7).Value = 1.604000333333333
$ws.Cells.Item(9, 8).Value = 4.812001
$ws.Cells.Item(9, 9).Value = 0.7668859909146659
$ws.Cells.Item(9, 10).Value = 0.7668859909146659
$ws.Cells.Item(9, 15).Value = 0.09956205276773258
$ws.Cells.Item(9, 16).Value = 0.09956205276773258
$ws.Cells.Item(9, 17).Value = 0.6132323181047779
$ws.Cells.Item(9, 18).Value = 5.519090862943001
$ws.Cells.Item(9, 19).Value = 0.07635274349428085
$ws.Cells.Item(9, 20).Value = 0.07635274349428085
$ws.Cells.Item(10, 7).Value = 1.604000333333333
$ws.Cells.Item(10, 8).Value = 4.812001
$ws.Cells.Item(10, 9).Value = 0.7668859909146659
$ws.Cells.Item(10, 10).Value = 0.7668859909146659
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 0.6666666666666666
$ws.Cells.Item(10, 13).Value = 0.655462
$ws.Cells.Item(10, 14).Value = 1.966386
$ws.Cells.Item(10, 15).Value = 0.1706949924222307
$ws.Cells.Item(10, 16).Value = 0.1706949924222307
$ws.Cells.Item(10, 17).Value = 1.051361266487333
$ws.Cells.Item(10, 18).Value = 9.462251398386
$ws.Cells.Item(10, 19).Value = 0.1309035984078938
$ws.Cells.Item(10, 20).Value = 0.1309035984078938
$ws.Cells.Item(11, 7).Value = 1.604000333333333
$ws.Cells.Item(11, 8).Value = 4.812001
$ws.Cells.Item(11, 9).Value = 0.7668859909146659
$ws.Cells.Item(11, 10).Value = 0.7668859909146659
$ws.Cells.Item(11, 13).Value = 1.828665666666667
$ws.Cells.Item(11, 14).Value = 5.485997
$ws.Cells.Item(11, 15).Value = 0.476219936646915
$ws.Cells.Item(11, 16).Value = 0.476219936646915
$ws.Cells.Item(11, 17).Value = 2.933180338888556
$ws.Cells.Item(11, 18).Value = 26.398623049997
$ws.Cells.Item(11, 19).Value = 0.3652063980087888
$ws.Cells.Item(11, 20).Value = 0.3652063980087888
$ws.Cells.Item(12, 7).Value = 1.604000333333333
$ws.Cells.Item(12, 8).Value = 4.812001
$ws.Cells.Item(12, 9).Value = 0.7668859909146659
$ws.Cells.Item(12, 10).Value = 0.7668859909146659
$ws.Cells.Item(12, 11).Value = 1
$ws.Cells.Item(12, 12).Value = 0.3333333333333333
$ws.Cells.Item(12, 13).Value = 0.07852066666666667
$ws.Cells.Item(12, 14).Value = 0.235562
$ws.Cells.Item(12, 15).Value = 0.02044830150589229
$ws.Cells.Item(12, 16).Value = 0.02044830150589229
$ws.Cells.Item(12, 17).Value = 0.1259471755068889
$ws.Cells.Item(12, 18).Value = 1.133524579562
$ws.Cells.Item(12, 19).Value = 0.01568151596286807
$ws.Cells.Item(12, 20).Value = 0.01568151596286807
$ws.Cells.Item(13, 7).Value = 1.604000333333333
$ws.Cells.Item(13, 8).Value = 4.812001
$ws.Cells.Item(13, 9).Value = 0.7668859909146659
$ws.Cells.Item(13, 10).Value = 0.7668859909146659
$ws.Cells.Item(13, 13).Value = 0.7251310000000001
$ws.Cells.Item(13, 14).Value = 2.175393
$ws.Cells.Item(13, 15).Value = 0.1888381485885141
$ws.Cells.Item(13, 16).Value = 0.1888381485885141
$ws.Cells.Item(13, 17).Value = 1.163110365710334
$ws.Cells.Item(13, 18).Value = 10.467993291393
$ws.Cells.Item(13, 19).Value = 0.1448173307027935
$ws.Cells.Item(13, 20).Value = 0.1448173307027935
